# Excel COM-interop script applying the edit described by the diff:
#  - Update row 34: D34 (Fecha) 44491 -> 44516, J34 (Volumen) 100 -> 120
#  - Insert a new row 35 that is a copy of the original (pre-edit) row 34
#    (same Fecha=44491 and Volumen=100 as the old row 34, all other columns identical)
#  - Worksheet dimension grows from A1:R34 to A1:R35 automatically

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: append the new row 35, a duplicate of the current row 34 ---
# (do this first, while row 34 still holds its original values)
# Note: use .Value2 (not .Value) - this engine's .Value getter does not
# reliably surface the underlying scalar, while .Value2 round-trips cleanly.
$ws.Range("A35").Value2 = $ws.Range("A34").Value2
$ws.Range("B35").Value2 = $ws.Range("B34").Value2
$ws.Range("C35").Value2 = $ws.Range("C34").Value2
$ws.Range("D35").Value2 = $ws.Range("D34").Value2
$ws.Range("D35").NumberFormat = $ws.Range("D34").NumberFormat
$ws.Range("E35").Value2 = $ws.Range("E34").Value2
$ws.Range("F35").Value2 = $ws.Range("F34").Value2
$ws.Range("G35").Value2 = $ws.Range("G34").Value2
$ws.Range("H35").Value2 = $ws.Range("H34").Value2
$ws.Range("I35").Value2 = $ws.Range("I34").Value2
$ws.Range("J35").Value2 = $ws.Range("J34").Value2
$ws.Range("K35").Value2 = $ws.Range("K34").Value2
$ws.Range("L35").Value2 = $ws.Range("L34").Value2
$ws.Range("M35").Value2 = $ws.Range("M34").Value2
$ws.Range("N35").Value2 = $ws.Range("N34").Value2
$ws.Range("O35").Value2 = $ws.Range("O34").Value2
$ws.Range("P35").Value2 = $ws.Range("P34").Value2
$ws.Range("Q35").Value2 = $ws.Range("Q34").Value2
$ws.Range("R35").Value2 = $ws.Range("R34").Value2

# --- Step 2: update row 34 in place with the new Fecha / Volumen values ---
$ws.Range("D34").Value2 = 44516
$ws.Range("J34").Value2 = 120

Write-Host "Row 35 added and row 34 updated."
